# 📊 Horarios actualizados Línea 141 - 891
# Refresh of the scraped bus-arrival schedule (new scrape at 08:29:58):
# existing upcoming rows get revised ETAs/minutes, and newly scraped
# arrivals are appended at the bottom of each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- LP1912 ---
$ws1.Cells.Item(2, 1).Value = "Última actualización: 08:29:58"
$ws1.Cells.Item(3, 1).Value = "Total filas: 118"
$ws1.Cells.Item(9, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(10, 3).Value = "15_ABASTO"
$ws1.Cells.Item(38, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(39, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(45, 1).Value = "06:46:37"
$ws1.Cells.Item(45, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(45, 4).Value = 4
$ws1.Cells.Item(46, 1).Value = "05:16:02"
$ws1.Cells.Item(46, 3).Value = "17_ROMERO"
$ws1.Cells.Item(46, 4).Value = 94
$ws1.Cells.Item(65, 1).Value = "06:46:37"
$ws1.Cells.Item(65, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(65, 4).Value = 50
$ws1.Cells.Item(66, 1).Value = "07:12:47"
$ws1.Cells.Item(66, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(66, 4).Value = 24
$ws1.Cells.Item(75, 1).Value = "06:53:56"
$ws1.Cells.Item(75, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(75, 4).Value = 66
$ws1.Cells.Item(76, 1).Value = "07:50:33"
$ws1.Cells.Item(76, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(76, 4).Value = 9
$ws1.Cells.Item(85, 1).Value = "08:29:58"
$ws1.Cells.Item(85, 4).Value = 0
$ws1.Cells.Item(86, 1).Value = "08:29:58"
$ws1.Cells.Item(86, 4).Value = 4
$ws1.Cells.Item(89, 1).Value = "08:29:58"
$ws1.Cells.Item(89, 4).Value = 12
$ws1.Cells.Item(90, 1).Value = "08:29:58"
$ws1.Cells.Item(90, 2).Value = "08:45"
$ws1.Cells.Item(90, 3).Value = "10_OLMOS"
$ws1.Cells.Item(90, 4).Value = 16
$ws1.Cells.Item(91, 1).Value = "08:29:58"
$ws1.Cells.Item(91, 2).Value = "08:47"
$ws1.Cells.Item(91, 4).Value = 18
$ws1.Cells.Item(93, 2).Value = "08:48"
$ws1.Cells.Item(93, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(93, 4).Value = 38
$ws1.Cells.Item(94, 1).Value = "08:29:58"
$ws1.Cells.Item(94, 2).Value = "08:50"
$ws1.Cells.Item(94, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(94, 4).Value = 21
$ws1.Cells.Item(95, 2).Value = "08:51"
$ws1.Cells.Item(95, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(95, 4).Value = 41
$ws1.Cells.Item(96, 1).Value = "08:29:58"
$ws1.Cells.Item(96, 2).Value = "08:59"
$ws1.Cells.Item(96, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(96, 4).Value = 30
$ws1.Cells.Item(97, 2).Value = "09:00"
$ws1.Cells.Item(97, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(97, 4).Value = 50
$ws1.Cells.Item(98, 1).Value = "08:29:58"
$ws1.Cells.Item(98, 2).Value = "09:01"
$ws1.Cells.Item(98, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(98, 4).Value = 32
$ws1.Cells.Item(99, 1).Value = "07:38:30"
$ws1.Cells.Item(99, 2).Value = "09:02"
$ws1.Cells.Item(99, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(99, 4).Value = 84
$ws1.Cells.Item(100, 1).Value = "08:29:58"
$ws1.Cells.Item(100, 2).Value = "09:03"
$ws1.Cells.Item(100, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(100, 4).Value = 34
$ws1.Cells.Item(101, 1).Value = "08:29:58"
$ws1.Cells.Item(101, 2).Value = "09:03"
$ws1.Cells.Item(101, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(101, 4).Value = 34
$ws1.Cells.Item(102, 1).Value = "08:10:38"
$ws1.Cells.Item(102, 2).Value = "09:10"
$ws1.Cells.Item(102, 4).Value = 60
$ws1.Cells.Item(103, 1).Value = "07:50:33"
$ws1.Cells.Item(103, 2).Value = "09:12"
$ws1.Cells.Item(103, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(103, 4).Value = 82
$ws1.Cells.Item(104, 1).Value = "08:29:58"
$ws1.Cells.Item(104, 2).Value = "09:14"
$ws1.Cells.Item(104, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(104, 4).Value = 45
$ws1.Cells.Item(105, 1).Value = "07:38:30"
$ws1.Cells.Item(105, 2).Value = "09:15"
$ws1.Cells.Item(105, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(105, 4).Value = 97
$ws1.Cells.Item(106, 1).Value = "08:29:58"
$ws1.Cells.Item(106, 2).Value = "09:16"
$ws1.Cells.Item(106, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(106, 4).Value = 47
$ws1.Cells.Item(107, 1).Value = "08:29:58"
$ws1.Cells.Item(107, 2).Value = "09:18"
$ws1.Cells.Item(107, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(107, 4).Value = 49
$ws1.Cells.Item(108, 1).Value = "07:50:33"
$ws1.Cells.Item(108, 2).Value = "09:19"
$ws1.Cells.Item(108, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(108, 4).Value = 89
$ws1.Cells.Item(109, 1).Value = "08:29:58"
$ws1.Cells.Item(109, 2).Value = "09:26"
$ws1.Cells.Item(109, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(109, 4).Value = 57
$ws1.Cells.Item(110, 1).Value = "08:29:58"
$ws1.Cells.Item(110, 2).Value = "09:28"
$ws1.Cells.Item(110, 4).Value = 59
$ws1.Cells.Item(111, 2).Value = "09:29"
$ws1.Cells.Item(111, 3).Value = "10_OLMOS"
$ws1.Cells.Item(111, 4).Value = 79
$ws1.Cells.Item(112, 1).Value = "08:29:58"
$ws1.Cells.Item(112, 2).Value = "09:33"
$ws1.Cells.Item(112, 3).Value = "15_ABASTO"
$ws1.Cells.Item(112, 4).Value = 64
$ws1.Cells.Item(113, 1).Value = "08:10:38"
$ws1.Cells.Item(113, 2).Value = "09:34"
$ws1.Cells.Item(113, 3).Value = "15_ABASTO"
$ws1.Cells.Item(113, 4).Value = 84
$ws1.Cells.Item(113, 5).Value = "LP1912"
$ws1.Cells.Item(114, 1).Value = "08:29:58"
$ws1.Cells.Item(114, 2).Value = "09:44"
$ws1.Cells.Item(114, 3).Value = "14_ABASTO"
$ws1.Cells.Item(114, 4).Value = 75
$ws1.Cells.Item(114, 5).Value = "LP1912"
$ws1.Cells.Item(115, 1).Value = "08:29:58"
$ws1.Cells.Item(115, 2).Value = "09:48"
$ws1.Cells.Item(115, 3).Value = "15_ABASTO"
$ws1.Cells.Item(115, 4).Value = 79
$ws1.Cells.Item(115, 5).Value = "LP1912"
$ws1.Cells.Item(116, 1).Value = "08:10:38"
$ws1.Cells.Item(116, 2).Value = "09:49"
$ws1.Cells.Item(116, 3).Value = "15_ABASTO"
$ws1.Cells.Item(116, 4).Value = 99
$ws1.Cells.Item(116, 5).Value = "LP1912"
$ws1.Cells.Item(117, 1).Value = "08:29:58"
$ws1.Cells.Item(117, 2).Value = "09:50"
$ws1.Cells.Item(117, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(117, 4).Value = 81
$ws1.Cells.Item(117, 5).Value = "LP1912"
$ws1.Cells.Item(118, 1).Value = "08:10:38"
$ws1.Cells.Item(118, 2).Value = "09:51"
$ws1.Cells.Item(118, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(118, 4).Value = 101
$ws1.Cells.Item(118, 5).Value = "LP1912"
$ws1.Cells.Item(119, 1).Value = "08:29:58"
$ws1.Cells.Item(119, 2).Value = "09:56"
$ws1.Cells.Item(119, 3).Value = "10_OLMOS"
$ws1.Cells.Item(119, 4).Value = 87
$ws1.Cells.Item(119, 5).Value = "LP1912"
$ws1.Cells.Item(120, 1).Value = "08:29:58"
$ws1.Cells.Item(120, 2).Value = "10:03"
$ws1.Cells.Item(120, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(120, 4).Value = 94
$ws1.Cells.Item(120, 5).Value = "LP1912"
$ws1.Cells.Item(121, 1).Value = "08:29:58"
$ws1.Cells.Item(121, 2).Value = "10:08"
$ws1.Cells.Item(121, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(121, 4).Value = 99
$ws1.Cells.Item(121, 5).Value = "LP1912"
$ws1.Cells.Item(122, 1).Value = "08:29:58"
$ws1.Cells.Item(122, 2).Value = "10:18"
$ws1.Cells.Item(122, 3).Value = "17_ROMERO"
$ws1.Cells.Item(122, 4).Value = 109
$ws1.Cells.Item(122, 5).Value = "LP1912"
$ws1.Cells.Item(123, 1).Value = "08:29:58"
$ws1.Cells.Item(123, 2).Value = "10:20"
$ws1.Cells.Item(123, 3).Value = "10_OLMOS"
$ws1.Cells.Item(123, 4).Value = 111
$ws1.Cells.Item(123, 5).Value = "LP1912"

# --- LP1912-215 ---
$ws2.Cells.Item(2, 1).Value = "Última actualización: 08:29:58"
$ws2.Cells.Item(18, 1).Value = "08:29:58"
$ws2.Cells.Item(18, 4).Value = 4
$ws2.Cells.Item(20, 1).Value = "08:29:58"
$ws2.Cells.Item(20, 4).Value = 18
$ws2.Cells.Item(22, 1).Value = "08:29:58"
$ws2.Cells.Item(22, 4).Value = 30
$ws2.Cells.Item(23, 1).Value = "08:29:58"
$ws2.Cells.Item(23, 4).Value = 49
$ws2.Cells.Item(25, 1).Value = "08:29:58"
$ws2.Cells.Item(25, 4).Value = 94

# --- 6203-6173 ---
$ws3.Cells.Item(2, 1).Value = "Última actualización: 08:29:58"
$ws3.Cells.Item(3, 1).Value = "Total filas: 13"
$ws3.Cells.Item(14, 1).Value = "08:29:58"
$ws3.Cells.Item(14, 4).Value = 22
$ws3.Cells.Item(16, 1).Value = "08:29:58"
$ws3.Cells.Item(16, 4).Value = 86
$ws3.Cells.Item(17, 1).Value = "08:29:58"
$ws3.Cells.Item(17, 2).Value = "10:10"
$ws3.Cells.Item(17, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(17, 4).Value = 101
$ws3.Cells.Item(17, 5).Value = "L6173"
$ws3.Cells.Item(18, 1).Value = "08:29:58"
$ws3.Cells.Item(18, 2).Value = "10:21"
$ws3.Cells.Item(18, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(18, 4).Value = 112
$ws3.Cells.Item(18, 5).Value = "L6173"
